# NATMI re-run (per Dr Hou's advice): adds the "M2" sending cluster and
# refreshes all ligand/receptor specificity metrics for Efnb1-Erbb2,
# growing the results table from 9 to 12 data rows (A1:T10 -> A1:T13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb1"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 11.43712066666667
$ws.Range("H2").Value = 34.311362
$ws.Range("I2").Value = 0.5796330080444665
$ws.Range("J2").Value = 0.5796330080444665
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.155977333333333
$ws.Range("N2").Value = 9.467932
$ws.Range("O2").Value = 0.3579027849973545
$ws.Range("P2").Value = 0.3579027849973545
$ws.Range("Q2").Value = 36.09529358259822
$ws.Range("R2").Value = 324.857642243384
$ws.Range("S2").Value = 0.2074522678555085
$ws.Range("T2").Value = 0.2074522678555085

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb1"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 11.43712066666667
$ws.Range("H3").Value = 34.311362
$ws.Range("I3").Value = 0.5796330080444665
$ws.Range("J3").Value = 0.5796330080444665
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.165953666666667
$ws.Range("N3").Value = 9.497861
$ws.Range("O3").Value = 0.359034148472735
$ws.Range("P3").Value = 0.359034148472735
$ws.Range("Q3").Value = 36.20939411074244
$ws.Range("R3").Value = 325.884546996682
$ws.Range("S3").Value = 0.208108043469935
$ws.Range("T3").Value = 0.208108043469935

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efnb1"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 11.43712066666667
$ws.Range("H4").Value = 34.311362
$ws.Range("I4").Value = 0.5796330080444665
$ws.Range("J4").Value = 0.5796330080444665
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.496042666666666
$ws.Range("N4").Value = 7.488128
$ws.Range("O4").Value = 0.2830630665299106
$ws.Range("P4").Value = 0.2830630665299106
$ws.Range("Q4").Value = 28.54754116781511
$ws.Range("R4").Value = 256.927870510336
$ws.Range("S4").Value = 0.1640726967190231
$ws.Range("T4").Value = 0.1640726967190231

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb1"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.754308333333334
$ws.Range("H5").Value = 17.262925
$ws.Range("I5").Value = 0.2916282118266253
$ws.Range("J5").Value = 0.2916282118266253
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.155977333333333
$ws.Range("N5").Value = 9.467932
$ws.Range("O5").Value = 0.3579027849973545
$ws.Range("P5").Value = 0.3579027849973545
$ws.Range("Q5").Value = 18.16046666901111
$ws.Range("R5").Value = 163.4442000211
$ws.Range("S5").Value = 0.1043745491965476
$ws.Range("T5").Value = 0.1043745491965476

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efnb1"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.754308333333334
$ws.Range("H6").Value = 17.262925
$ws.Range("I6").Value = 0.2916282118266253
$ws.Range("J6").Value = 0.2916282118266253
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.165953666666667
$ws.Range("N6").Value = 9.497861
$ws.Range("O6").Value = 0.359034148472735
$ws.Range("P6").Value = 0.359034148472735
$ws.Range("Q6").Value = 18.21787356704722
$ws.Range("R6").Value = 163.960862103425
$ws.Range("S6").Value = 0.1047044867037988
$ws.Range("T6").Value = 0.1047044867037988

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efnb1"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.754308333333334
$ws.Range("H7").Value = 17.262925
$ws.Range("I7").Value = 0.2916282118266253
$ws.Range("J7").Value = 0.2916282118266253
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.496042666666666
$ws.Range("N7").Value = 7.488128
$ws.Range("O7").Value = 0.2830630665299106
$ws.Range("P7").Value = 0.2830630665299106
$ws.Range("Q7").Value = 14.36299911715556
$ws.Range("R7").Value = 129.2669920544
$ws.Range("S7").Value = 0.0825491759262789
$ws.Range("T7").Value = 0.0825491759262789

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Efnb1"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.2670106666666667
$ws.Range("H8").Value = 0.8010320000000001
$ws.Range("I8").Value = 0.0135320943453039
$ws.Range("J8").Value = 0.0135320943453039
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.155977333333333
$ws.Range("N8").Value = 9.467932
$ws.Range("O8").Value = 0.3579027849973545
$ws.Range("P8").Value = 0.3579027849973545
$ws.Range("Q8").Value = 0.8426796117582221
$ws.Range("R8").Value = 7.584116505824
$ws.Range("S8").Value = 0.004843174253031217
$ws.Range("T8").Value = 0.004843174253031218

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Efnb1"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.2670106666666667
$ws.Range("H9").Value = 0.8010320000000001
$ws.Range("I9").Value = 0.0135320943453039
$ws.Range("J9").Value = 0.0135320943453039
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.165953666666667
$ws.Range("N9").Value = 9.497861
$ws.Range("O9").Value = 0.359034148472735
$ws.Range("P9").Value = 0.359034148472735
$ws.Range("Q9").Value = 0.8453433991724445
$ws.Range("R9").Value = 7.608090592552001
$ws.Range("S9").Value = 0.004858483970318897
$ws.Range("T9").Value = 0.004858483970318898

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Efnb1"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.2670106666666667
$ws.Range("H10").Value = 0.8010320000000001
$ws.Range("I10").Value = 0.0135320943453039
$ws.Range("J10").Value = 0.0135320943453039
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.496042666666666
$ws.Range("N10").Value = 7.488128
$ws.Range("O10").Value = 0.2830630665299106
$ws.Range("P10").Value = 0.2830630665299106
$ws.Range("Q10").Value = 0.6664700164551111
$ws.Range("R10").Value = 5.998230148096001
$ws.Range("S10").Value = 0.003830436121953785
$ws.Range("T10").Value = 0.003830436121953786

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Efnb1"
$ws.Range("C11").Value = "Erbb2"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.273219
$ws.Range("H11").Value = 6.819656999999999
$ws.Range("I11").Value = 0.1152066857836043
$ws.Range("J11").Value = 0.1152066857836043
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.155977333333333
$ws.Range("N11").Value = 9.467932
$ws.Range("O11").Value = 0.3579027849973545
$ws.Range("P11").Value = 0.3579027849973545
$ws.Range("Q11").Value = 7.174227637702665
$ws.Range("R11").Value = 64.56804873932398
$ws.Range("S11").Value = 0.0412327936922671
$ws.Range("T11").Value = 0.04123279369226711

$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Efnb1"
$ws.Range("C12").Value = "Erbb2"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.273219
$ws.Range("H12").Value = 6.819656999999999
$ws.Range("I12").Value = 0.1152066857836043
$ws.Range("J12").Value = 0.1152066857836043
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.165953666666667
$ws.Range("N12").Value = 9.497861
$ws.Range("O12").Value = 0.359034148472735
$ws.Range("P12").Value = 0.359034148472735
$ws.Range("Q12").Value = 7.196906028186333
$ws.Range("R12").Value = 64.772154253677
$ws.Range("S12").Value = 0.04136313432868231
$ws.Range("T12").Value = 0.04136313432868232

$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Efnb1"
$ws.Range("C13").Value = "Erbb2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.273219
$ws.Range("H13").Value = 6.819656999999999
$ws.Range("I13").Value = 0.1152066857836043
$ws.Range("J13").Value = 0.1152066857836043
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.496042666666666
$ws.Range("N13").Value = 7.488128
$ws.Range("O13").Value = 0.2830630665299106
$ws.Range("P13").Value = 0.2830630665299106
$ws.Range("Q13").Value = 5.674051614677332
$ws.Range("R13").Value = 51.06646453209599
$ws.Range("S13").Value = 0.03261075776265489
$ws.Range("T13").Value = 0.0326107577626549
